$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting first for Price cells whose new value would
# otherwise be auto-detected as a Number by Excel (losing formatting
# such as trailing zeros), matching the original inlineStr text cells.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns with the latest values.
$ws.Range("D2").Value = "29.904.85"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.873.77"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "0.7434"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "241.92"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").Value = "0.9980"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "0.3153"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.07175"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "24.79"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").Value = "0.08429"
$ws.Range("E11").Value = "  -5.57%  "
$ws.Range("D12").Value = "0.7544"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "5.438"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "1.870.82"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "92.75"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "29.908.45"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "6.059"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "13.61"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "244.13"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "0.000007830"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "0.9977"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "2.112.65"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "7.989"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("D24").Value = "0.9988"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "0.1572"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "9.306"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "164.20"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "2.031"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("D31").Value = "4.610"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "4.270"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "0.05326"
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "0.7555"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").Value = "0.9969"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "2.693"
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "2.745"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "0.4485"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "1.103.47"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "6.106"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "72.40"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").Value = "0.8603"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "0.9989"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "103.23"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "7.697"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").Value = "3.049"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "2.012.45"
$ws.Range("E51").Value = "  -0.77%  "
